# Weekly CompStat refresh: new volume/date header + updated 123rd Precinct crime counts.
function Set-NumericCell($ws, $cell, $val, $styleRef) {
    # Pull the numeric-style formatting (e.g. "#,##0") from a donor cell, then write the number.
    $ws.Range($styleRef).Copy() | Out-Null
    $ws.Range($cell).PasteSpecial(-4122) | Out-Null
    $ws.Range($cell).Value = $val
}

function Set-TextCell($ws, $cell, $val, $styleRef) {
    # Force text (not auto-coerced-to-number) semantics, then restore the donor text style.
    $ws.Range($cell).NumberFormat = "@"
    $ws.Range($cell).Value = $val
    $ws.Range($styleRef).Copy() | Out-Null
    $ws.Range($cell).PasteSpecial(-4122) | Out-Null
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Masthead: volume number + week-covered date range ---
$ws.Range("A8").Value = "Volume 32   Number  37"
$ws.Range("C9").Value = "Report Covering the Week  9/8/2025  Through  9/14/2025"

# --- Row 15: Rape ---
Set-NumericCell $ws "D15" 1 "I15"
Set-NumericCell $ws "E15" -100 "K15"
Set-NumericCell $ws "G15" 1 "I15"
Set-NumericCell $ws "H15" -100 "K15"
$ws.Range("J15").Value = 5
$ws.Range("K15").Value = 0
# --- Row 16: Robbery ---
Set-NumericCell $ws "D16" 1 "I15"
Set-NumericCell $ws "E16" -100 "K15"
Set-NumericCell $ws "G16" 1 "I15"
Set-NumericCell $ws "H16" 0 "K15"
$ws.Range("J16").Value = 19
$ws.Range("K16").Value = 0
$ws.Range("N16").Value = -54.761904761904
# --- Row 17: Fel. Assault ---
$ws.Range("C17").Value = 1
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 9
$ws.Range("G17").Value = 7
$ws.Range("H17").Value = 28.571428571428
$ws.Range("I17").Value = 75
$ws.Range("J17").Value = 66
$ws.Range("K17").Value = 13.636363636363
$ws.Range("L17").Value = 15.384615384615
$ws.Range("M17").Value = 134.375
$ws.Range("N17").Value = -7.407407407407
# --- Row 18: Burglary ---
Set-TextCell $ws "C18" "0" "C14"
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 4
$ws.Range("H18").Value = -20
$ws.Range("I18").Value = 36
$ws.Range("J18").Value = 26
$ws.Range("K18").Value = 38.461538461538
$ws.Range("L18").Value = -10
$ws.Range("M18").Value = -52.631578947368
$ws.Range("N18").Value = -83.48623853211
# --- Row 19: Gr. Larceny ---
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 5
$ws.Range("E19").Value = 40
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = 75
$ws.Range("I19").Value = 211
$ws.Range("J19").Value = 208
$ws.Range("K19").Value = 1.442307692307
$ws.Range("L19").Value = 4.455445544554
$ws.Range("M19").Value = 113.131313131313
$ws.Range("N19").Value = 33.544303797468
# --- Row 20: G.L.A. ---
$ws.Range("G20").Value = 5
$ws.Range("J20").Value = 30
$ws.Range("K20").Value = -43.333333333333
$ws.Range("L20").Value = -70.689655172413
$ws.Range("N20").Value = -96.633663366336
# --- Row 21: TOTAL ---
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = -20
$ws.Range("F21").Value = 49
$ws.Range("G21").Value = 39
$ws.Range("H21").Value = 25.641025641025
$ws.Range("I21").Value = 363
$ws.Range("J21").Value = 356
$ws.Range("K21").Value = 1.966292134831
$ws.Range("L21").Value = -4.724409448818
$ws.Range("M21").Value = 43.478260869565
$ws.Range("N21").Value = -63.988095238095
# --- Row 24: Petit Larceny ---
$ws.Range("C24").Value = 7
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 40
$ws.Range("F24").Value = 39
$ws.Range("G24").Value = 24
$ws.Range("H24").Value = 62.5
$ws.Range("I24").Value = 307
$ws.Range("J24").Value = 255
$ws.Range("K24").Value = 20.392156862745
$ws.Range("L24").Value = -8.358208955223
$ws.Range("M24").Value = -20.25974025974
# --- Row 25: Retail Theft ---
$ws.Range("C25").Value = 2
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 16
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = 77.777777777777
$ws.Range("I25").Value = 180
$ws.Range("J25").Value = 116
$ws.Range("K25").Value = 55.172413793103
$ws.Range("L25").Value = 21.621621621621
# --- Row 26: Misd. Assault ---
$ws.Range("F26").Value = 18
$ws.Range("G26").Value = 10
$ws.Range("H26").Value = 80
$ws.Range("I26").Value = 147
$ws.Range("J26").Value = 113
$ws.Range("K26").Value = 30.088495575221
$ws.Range("L26").Value = 2.797202797202
$ws.Range("M26").Value = -3.28947368421
# --- Row 27: UCR Rape* ---
Set-NumericCell $ws "D27" 1 "I15"
Set-NumericCell $ws "E27" -100 "K15"
Set-NumericCell $ws "G27" 1 "I15"
Set-NumericCell $ws "H27" 0 "K15"
$ws.Range("J27").Value = 9
$ws.Range("K27").Value = -33.333333333333
# --- Row 28: Other Sex Crimes ---
Set-NumericCell $ws "C28" 1 "I15"
Set-NumericCell $ws "D28" 1 "I15"
Set-NumericCell $ws "E28" 0 "K15"
$ws.Range("F28").Value = 3
Set-NumericCell $ws "G28" 1 "I15"
Set-NumericCell $ws "H28" 200 "K15"
$ws.Range("I28").Value = 12
$ws.Range("J28").Value = 10
$ws.Range("K28").Value = 20
$ws.Range("L28").Value = -7.692307692307
# --- Row 33: Hate Crimes ---
Set-TextCell $ws "D33" "0" "C14"
Set-TextCell $ws "E33" "***.*" "C14"

$excel.CutCopyMode = 0
